$d = $word.ActiveDocument

# --- Edit 1 (SVM paragraph): "for both dataset," -> "for both datasets," ---
# (net text change; the underlying XML in the target also re-splits runs and
#  shifts the grammar-check bracket, but the visible text is simply the
#  insertion of the missing "s".)
$null = $d.Content.Find.Execute(
    "for both dataset, the optimal kernel",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for both datasets, the optimal kernel", 2)

# --- Edit 2 (KNN paragraph): "less than six" -> "less than 6" ---
$null = $d.Content.Find.Execute(
    "less than six indicating overfitting",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "less than 6 indicating overfitting", 2)

# --- Edit 3: append a new "Neural Network" section (three new paragraphs)
#     right after the KNN learning-curve paragraph, before the trailing
#     empty paragraph that precedes the section break. ---
$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n - 1)

$null = $target.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Item($n)
$pHeading.Range.Text = "Neural Network"

$null = $pHeading.Range.InsertParagraphAfter()
$pBody1 = $d.Paragraphs.Item($n + 1)
$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$body1 = "We will be using the default activation function " + $quoteOpen + "ReLU" + $quoteClose + " for our Neural Network Learner. During our early model exploration, both datasets benefit when around learning rate is around 0.01. When the learning rate is set to low, the model underfit for both datasets. However, as learning rate increases, underfitting decreases and test accuracy increases up until 0.1. The test accuracy starts decreasing after learning rate goes beyond 0.1. Similarly, higher the number of hidden layers better the accuracy is for both datasets. We will now use GridSearchCV to find the most optimal learning rate and optimal hidden layers."
$pBody1.Range.Text = $body1

$null = $pBody1.Range.InsertParagraphAfter()
$pBody2 = $d.Paragraphs.Item($n + 2)
$body2 = "The learning curve of dataset 1 suggests that due to high variance of the data, adding more training instances will not necessarily increase the cross-validation score.  However, adding more test instances does remove bias and variance in dataset 1. In dataset 2, the cross-validation score increases as more training instances are introduced indicating that the model will benefit from adding more training data. "
$pBody2.Range.Text = $body2
